$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.726.25"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "1.574.31"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'213.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'44.66"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "'24.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.799.98"
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "1.573.82"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "28.727.22"
$ws.Range("E16").Value = "  +1.76%  "
$ws.Range("D17").Value = "'3.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "'62.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "'230.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("E23").Value = "  -4.56%  "
$ws.Range("D24").Value = "'9.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "'2.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.65%  "
$ws.Range("D26").Value = "'151.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'15.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").Value = "'6.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D35").Value = "1.397.36"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "'1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("E39").Value = "  +3.01%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  -2.60%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'0.795"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").Value = "'0.0469"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "'0.961"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("D48").Value = "'63.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").Value = "1.711.26"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "'86.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  +0.21%  "
